# Weekly fruit/vegetable data update:
# Insert a new daily record (2021-10-05, date serial 44474) for "Sandia"
# (Extra/Primera quality) sourced from "Perú", as a new row right above
# the previously-first 2021-03-31 batch of records (old row 182).
# All subsequent rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 182; everything from 182 downward shifts to 183+
$ws.Rows("182").Insert()

# Populate the newly inserted row 182 with the new record
$ws.Range("A182").Value = 9
$ws.Range("B182").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C182").Value = "Metropolitana"
$ws.Range("D182").Value = 44474
$ws.Range("E182").Value = 13
$ws.Range("F182").Value = 100112028
$ws.Range("G182").Value = "Sandia"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 250
$ws.Range("K182").Value = 800
$ws.Range("L182").Value = 1000
$ws.Range("M182").Value = 900
$ws.Range("N182").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O182").Value = "Perú"
$ws.Range("P182").Value = 900
$ws.Range("Q182").Value = 1
$ws.Range("R182").Value = "Hortaliza"
